$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 17:10"

# Refresh case statistics for Estados Unidos (row 4), Alemania (row 11), Canada (row 16)
$ws.Range("B4").Value = 1749726
$ws.Range("C4").Value = 3923
$ws.Range("D4").Value = 490262
$ws.Range("E4").Value = 1157215
$ws.Range("G4").Value = 142
$ws.Range("H4").Value = 102249

$ws.Range("B11").Value = 182150
$ws.Range("C11").Value = 255
$ws.Range("E11").Value = 10402
$ws.Range("G11").Value = 15
$ws.Range("H11").Value = 8548

$ws.Range("B16").Value = 87902
$ws.Range("C16").Value = 383
$ws.Range("D16").Value = 46465
$ws.Range("E16").Value = 34638
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 6799

# Re-sort tied countries (equal "Casos totales") by updating country names and
# their accompanying stats so each row carries the right country's full data
$ws.Range("A197").Value = "Curazao"
$ws.Range("D197").Value = 14
$ws.Range("H197").Value = 1

$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Santa Lucia"

$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Nueva Caledonia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "San Bartolome"

$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
